$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2197.6
$ws.Range("I86").Value = 1833
$ws.Range("J86").Value = 2744.5
$ws.Range("K86").Value = 1833
$ws.Range("L86").Value = 2744.5
$ws.Range("M86").Value = -710
$ws.Range("N86").Value = -4990.5

$ws.Range("H88").Value = 872.5714
$ws.Range("I88").Value = 503
$ws.Range("J88").Value = 934.1667
$ws.Range("K88").Value = 503
$ws.Range("L88").Value = 934.1667
$ws.Range("M88").Value = -97
$ws.Range("N88").Value = -1746.1667

$ws.Range("H89").Value = 2197.6
$ws.Range("I89").Value = 1833
$ws.Range("J89").Value = 2744.5
$ws.Range("K89").Value = 9165
$ws.Range("L89").Value = 13722.5
$ws.Range("M89").Value = -3549
$ws.Range("N89").Value = -24954.5

$ws.Range("H91").Value = 872.5714
$ws.Range("I91").Value = 503
$ws.Range("J91").Value = 934.1667
$ws.Range("K91").Value = 503
$ws.Range("L91").Value = 934.1667
$ws.Range("M91").Value = 901
$ws.Range("N91").Value = -3742.1667

$ws.Range("H98").Value = 4522
$ws.Range("I98").Value = 4785.1816
$ws.Range("K98").Value = 4785.1816
$ws.Range("M98").Value = -3287.1816

$ws.Range("H99").Value = 252.25
$ws.Range("I99").Value = 252.25
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 756.75
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 741.25
$ws.Range("N99").ClearContents()

$ws.Range("H100").Value = 1379.6842
$ws.Range("I100").Value = 1512.2941
$ws.Range("J100").Value = 252.5
$ws.Range("K100").Value = 1512.2941
$ws.Range("L100").Value = 252.5
$ws.Range("M100").Value = -971.2941000000001
$ws.Range("N100").Value = -1334.5

$ws.Range("H107").Value = 783.5263
$ws.Range("I107").Value = 758.0714
$ws.Range("K107").Value = 758.0714
$ws.Range("M107").Value = 1161.9286

$ws.Range("H111").Value = 13421.412
$ws.Range("I111").Value = 12097.333
$ws.Range("J111").Value = 14911
$ws.Range("K111").Value = 36291.999
$ws.Range("L111").Value = 44733
$ws.Range("M111").Value = -33224.999
$ws.Range("N111").Value = -50867

$ws.Range("H112").Value = 25867.5
$ws.Range("I112").Value = 969.7143
$ws.Range("J112").Value = 30577.893
$ws.Range("K112").Value = 2909.1429
$ws.Range("L112").Value = 91733.679
$ws.Range("M112").Value = -1801.1429
$ws.Range("N112").Value = -93949.679

$ws.Range("H115").Value = 298.76923
$ws.Range("I115").Value = 298.76923
$ws.Range("K115").Value = 896.30769
$ws.Range("M115").Value = 670.69231

$ws.Range("H122").Value = 4522
$ws.Range("I122").Value = 4785.1816
$ws.Range("K122").Value = 14355.5448
$ws.Range("M122").Value = -11905.5448

$ws.Range("H141").Value = 959.9231
$ws.Range("I141").Value = 959.9231
$ws.Range("K141").Value = 2879.7693
$ws.Range("M141").Value = 2300.2307

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7926.1885
$ws.Range("J32").Value = 33018
$ws.Range("L32").Value = 33018
$ws.Range("N32").Value = -33592

$ws.Range("H61").Value = 4076.7827
$ws.Range("I61").Value = 3670.3333
$ws.Range("K61").Value = 3670.3333
$ws.Range("M61").Value = -3458.3333

$ws.Range("H74").Value = 7868.3335
$ws.Range("I74").Value = 1861.1578
$ws.Range("J74").Value = 18244.363
$ws.Range("K74").Value = 1861.1578
$ws.Range("L74").Value = 18244.363
$ws.Range("M74").Value = -987.1578
$ws.Range("N74").Value = -19992.363

$ws.Range("H77").Value = 7868.3335
$ws.Range("I77").Value = 1861.1578
$ws.Range("J77").Value = 18244.363
$ws.Range("K77").Value = 9305.789000000001
$ws.Range("L77").Value = 91221.815
$ws.Range("M77").Value = -4937.789000000001
$ws.Range("N77").Value = -99957.815

$ws.Range("H102").Value = 3089
$ws.Range("I102").Value = 2778.6
$ws.Range("J102").Value = 3433.889
$ws.Range("K102").Value = 2778.6
$ws.Range("L102").Value = 3433.889
$ws.Range("M102").Value = -1156.6
$ws.Range("N102").Value = -6677.889

$ws.Range("H136").Value = 4076.7827
$ws.Range("I136").Value = 3670.3333
$ws.Range("K136").Value = 11010.9999
$ws.Range("M136").Value = -8460.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3787.6592
$ws.Range("J20").Value = 4892
$ws.Range("L20").Value = 4892
$ws.Range("N20").Value = -5386

$ws.Range("H86").Value = 3351.7222
$ws.Range("I86").Value = 3173.923
$ws.Range("J86").Value = 3814
$ws.Range("K86").Value = 3173.923
$ws.Range("L86").Value = 3814
$ws.Range("M86").Value = -2050.923
$ws.Range("N86").Value = -6060

$ws.Range("H89").Value = 3351.7222
$ws.Range("I89").Value = 3173.923
$ws.Range("J89").Value = 3814
$ws.Range("K89").Value = 15869.615
$ws.Range("L89").Value = 19070
$ws.Range("M89").Value = -10253.615
$ws.Range("N89").Value = -30302

$ws.Range("H94").Value = 1351.6562
$ws.Range("I94").Value = 926.25
$ws.Range("J94").Value = 2627.875
$ws.Range("K94").Value = 926.25
$ws.Range("L94").Value = 2627.875
$ws.Range("M94").Value = -475.25
$ws.Range("N94").Value = -3529.875

$ws.Range("H99").Value = 4166.9
$ws.Range("I99").Value = 2026.4166
$ws.Range("J99").Value = 7377.625
$ws.Range("K99").Value = 2026.4166
$ws.Range("L99").Value = 7377.625
$ws.Range("M99").Value = -528.4166
$ws.Range("N99").Value = -10373.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H54").Value = 49750
$ws.Range("J54").Value = 49750
$ws.Range("L54").Value = 49750
$ws.Range("N54").Value = -51066

$ws.Range("H58").Value = 3462.862
$ws.Range("I58").Value = 3202.1
$ws.Range("K58").Value = 3202.1
$ws.Range("M58").Value = -2999.1

$ws.Range("H60").Value = 19687.375
$ws.Range("J60").Value = 19642.857
$ws.Range("L60").Value = 19642.857
$ws.Range("N60").Value = -20664.857

$ws.Range("H132").Value = 3716.1667
$ws.Range("I132").Value = 3764.0833
$ws.Range("K132").Value = 11292.2499
$ws.Range("M132").Value = -8762.249899999999

$ws.Range("H134").Value = 20185.04
$ws.Range("I134").Value = 11441.3
$ws.Range("K134").Value = 34323.89999999999
$ws.Range("M134").Value = -31788.89999999999

$ws.Range("H136").Value = 3462.862
$ws.Range("I136").Value = 3202.1
$ws.Range("K136").Value = 9606.299999999999
$ws.Range("M136").Value = -7056.299999999999

$ws.Range("H141").Value = 468065
$ws.Range("J141").Value = 733442
$ws.Range("L141").Value = 733442
$ws.Range("N141").Value = -743802

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H107").Value = 766.1667
$ws.Range("I107").Value = 839.4
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 2518.2
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = -598.1999999999998
$ws.Range("N107").Value = -5040

$ws.Range("H121").Value = 738
$ws.Range("J121").Value = 915
$ws.Range("L121").Value = 2745
$ws.Range("N121").Value = -5365

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 113.625
$ws.Range("I2").Value = 77.8421
$ws.Range("J2").Value = 249.6
$ws.Range("K2").Value = 77.8421
$ws.Range("L2").Value = 249.6
$ws.Range("M2").Value = 35.1579
$ws.Range("N2").Value = -475.6

$ws.Range("H70").Value = 4049.1538
$ws.Range("I70").Value = 3883.2
$ws.Range("K70").Value = 3883.2
$ws.Range("M70").Value = -3613.2

$ws.Range("H73").Value = 4049.1538
$ws.Range("I73").Value = 3883.2
$ws.Range("K73").Value = 3883.2
$ws.Range("M73").Value = -2947.2

$ws.Range("H80").Value = 2306.353
$ws.Range("I80").Value = 2588.625
$ws.Range("J80").Value = 2055.4443
$ws.Range("K80").Value = 2588.625
$ws.Range("L80").Value = 2055.4443
$ws.Range("M80").Value = -1590.625
$ws.Range("N80").Value = -4051.4443

$ws.Range("H83").Value = 2306.353
$ws.Range("I83").Value = 2588.625
$ws.Range("J83").Value = 2055.4443
$ws.Range("K83").Value = 12943.125
$ws.Range("L83").Value = 10277.2215
$ws.Range("M83").Value = -7951.125
$ws.Range("N83").Value = -20261.2215

$ws.Range("H107").Value = 461.25
$ws.Range("I107").Value = 456.75
$ws.Range("J107").Value = 474.75
$ws.Range("K107").Value = 456.75
$ws.Range("L107").Value = 474.75
$ws.Range("M107").Value = 1463.25
$ws.Range("N107").Value = -4314.75

$ws.Range("H126").Value = 17267.268
$ws.Range("I126").Value = 19307.309
$ws.Range("J126").Value = 4007
$ws.Range("K126").Value = 57921.927
$ws.Range("L126").Value = 12021
$ws.Range("M126").Value = -55451.927
$ws.Range("N126").Value = -16961

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4080.8125
$ws.Range("I136").Value = 3791.7693
$ws.Range("K136").Value = 11375.3079
$ws.Range("M136").Value = -8825.3079

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1151.75
$ws.Range("J96").Value = 1001.3333
$ws.Range("L96").Value = 1001.3333
$ws.Range("N96").Value = -3747.3333

$ws.Range("H136").Value = 3602.3572
$ws.Range("I136").Value = 3602.0908
$ws.Range("J136").Value = 3603.3333
$ws.Range("K136").Value = 10806.2724
$ws.Range("L136").Value = 10809.9999
$ws.Range("M136").Value = -8256.2724
$ws.Range("N136").Value = -15909.9999
